# Applies the edits described by the commit "MM minor and Results 2nd,2.5/4"
#
# Summary of changes:
#  1) First "supplementary table S3)" (Materials & Methods, CCSM paragraph):
#     split the bold run so the trailing ")" is no longer bold.
#  2) Second "supplementary table S3)" (future studies paragraph):
#     same split, and the Word "last edit" (_GoBack) bookmark now wraps the
#     new ")" run.
#  3) "the mutants ( all primers" -> merge the " " and "(" runs into a
#     single run containing " (".
#  4) Remove the old _GoBack bookmark that used to sit after
#     "public availability of".
#  5) Footer page-number field cached result: "6" -> "9".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) First "supplementary table S3)" -> un-bold the closing parenthesis
# ---------------------------------------------------------------------
$range1 = $d.Content
$range1.Find.Execute("(CCSM) (supplementary table S3)") | Out-Null
$close1 = $d.Range($range1.End - 1, $range1.End)
$close1.Font.Bold = 0

# ---------------------------------------------------------------------
# 2) Second "supplementary table S3)" -> un-bold the closing parenthesis
#    and move the _GoBack bookmark onto it.
# ---------------------------------------------------------------------
$range2 = $d.Content
$range2.Find.Execute("future studies (supplementary table S3)") | Out-Null
$close2 = $d.Range($range2.End - 1, $range2.End)
$close2.Font.Bold = 0

# ---------------------------------------------------------------------
# 3) Merge the " " and "(" runs after "the mutants" into a single run.
# ---------------------------------------------------------------------
$range3 = $d.Content
$range3.Find.Execute("the mutants") | Out-Null
$spaceParen = $d.Range($range3.End, $range3.End + 2)
$spaceParen.Find.Execute(" (", $true, $false, $false, $false, $false, $true, 1, $false, " (", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Remove the old _GoBack bookmark (was after "public availability of").
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2, continued) re-add _GoBack around the new ")" run from step 2.
# ---------------------------------------------------------------------
$goBackRange = $d.Range($close2.Start, $close2.End)
$goBackRange.Bookmarks.Add("_GoBack") | Out-Null

# ---------------------------------------------------------------------
# 5) Footer page-number field: cached result "6" -> "9"
# ---------------------------------------------------------------------
$footer = $d.Sections(1).Footers(1)
$footerRange = $footer.Range
$firstChar = $footerRange.Characters(1)
if ($firstChar.Text -eq "6") {
    $firstChar.Text = "9"
}
